# Apply bracket pick corrections on the "Bracket" worksheet.
# Each cell below holds the name of the animal advancing in the bracket;
# these are being corrected to the animal that actually should have advanced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$updates = @{
    "E12" = "Chequered elephant shrew";
    "M12" = "Palaeocastor";
    "D14" = "Chequered elephant shrew";
    "N14" = "Tent-making Bat";
    "G16" = "Okapi";
    "K16" = "Cathedral Termite";
    "E20" = "Fire-footed Rope Squirrel";
    "M20" = "Homo habilis";
    "F24" = "Numbat";
    "L24" = "Cathedral Termite";
    "D26" = "Numbat";
    "E28" = "Numbat";
    "M28" = "Cathedral Termite";
    "D30" = "Badger Bat";
    "D34" = "Sea Otter";
    "F40" = "Sea Otter";
    "D42" = "Silver Pika";
    "E44" = "Siberian Chipmunk";
    "M44" = "Siamang";
    "D46" = "Siberian Chipmunk";
    "N46" = "Siamang";
    "D50" = "Silky Anteater";
    "E52" = "Colo Colo Opossum";
    "L56" = "Greater Rhea";
    "D58" = "Bulldog Bat";
    "M60" = "Greater Rhea";
    "D62" = "Rock Hyrax ";
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
